# Remove the rows corresponding to Caso 915, 1068, 2676 and 5589.
# These are currently located at worksheet rows 2, 3, 5 and 9.
# Deleting from bottom to top keeps the remaining row numbers stable
# while we iterate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(9, 5, 3, 2)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
